$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Totais")

# New daily data to append: dates (as serials) Feb 15-22 2026 and their delivered totals
$dates = 46068,46069,46070,46071,46072,46073,46074,46075
$vals  = 383,328,370,336,349,318,268,310

$oldTotalRow = 16
$newTotalRow = $oldTotalRow + $dates.Count

# 1) Move the existing TOTAL row (currently row 16, label + formatted sum) down to its new
#    position (row 24), carrying over cell formatting (number format/font/border/fill).
$ws.Range("A" + $oldTotalRow + ":B" + $oldTotalRow).Copy($ws.Range("A" + $newTotalRow + ":B" + $newTotalRow))

# 2) Re-point the TOTAL formula to cover the full new data range before the old row 16 gets
#    overwritten with new daily data below.
$ws.Range("B" + $newTotalRow).Formula = "=SUM(B2:B" + ($newTotalRow - 1) + ")"

# 3) Fill rows 16-23 with the new daily data, copying the date column's format from the
#    last existing data row so the new dates render the same way.
for ($i = 0; $i -lt $dates.Count; $i++) {
    $r = $oldTotalRow + $i
    $ws.Range("A15").Copy($ws.Range("A" + $r))
    $ws.Range("A" + $r).Value = $dates[$i]
    $ws.Range("B" + $r).ClearFormats()
    $ws.Range("B" + $r).Value = $vals[$i]
}

# 4) Update the frozen pane / selection to match the new layout.
$activeWindow = $excel.ActiveWindow
$ws.Range("A5").Select()
$activeWindow.FreezePanes = $false
$activeWindow.FreezePanes = $true
$ws.Range("B25").Select()
